$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Config1")
$ws2 = $wb.Worksheets.Item("Config2")

# --- Config1 (sheet2): xalancbmk SRRIP/Hawkeye/OPTGen rows + zeusmp LRU row ---
$ws1.Range("C80").Value = 50000000
$ws1.Range("D80").Value = 241354593
$ws1.Range("E80").Value = 1587556
$ws1.Range("F80").Value = 162379
$ws1.Range("G80").Value = 1425177

$ws1.Range("C81").Value = 50000000
$ws1.Range("D81").Value = 228028533
$ws1.Range("E81").Value = 1587556
$ws1.Range("F81").Value = 363341
$ws1.Range("G81").Value = 1224215

$ws1.Range("C82").Value = 50000000
$ws1.Range("D82").Value = 228028533
$ws1.Range("E82").Value = 47204
$ws1.Range("F82").Value = 16412
$ws1.Range("G82").Formula = "=E82-F82"
$ws1.Range("J82").Formula = "=F82/E82"

$ws1.Range("C83").Value = 50000000
$ws1.Range("D83").Value = 43611395
$ws1.Range("E83").Value = 402709
$ws1.Range("F83").Value = 151183
$ws1.Range("G83").Value = 251526

# --- Config2 (sheet3): xalancbmk SRRIP/Hawkeye/OPTGen rows + zeusmp LRU row ---
$ws2.Range("C80").Value = 50000000
$ws2.Range("D80").Value = 192694662
$ws2.Range("E80").Value = 2168997
$ws2.Range("F80").Value = 161045
$ws2.Range("G80").Value = 2007952

$ws2.Range("C81").Value = 50000000
$ws2.Range("D81").Value = 186495428
$ws2.Range("E81").Value = 2168866
$ws2.Range("F81").Value = 334605
$ws2.Range("G81").Value = 1834261

$ws2.Range("C82").Value = 50000000
$ws2.Range("D82").Value = 186495428
$ws2.Range("E82").Value = 27979
$ws2.Range("F82").Value = 14942
$ws2.Range("G82").Formula = "=E82-F82"
$ws2.Range("J82").Formula = "=F82/E82"

$ws2.Range("C83").Value = 50000001
$ws2.Range("D83").Value = 34327145
$ws2.Range("E83").Value = 409177
$ws2.Range("F83").Value = 157503
$ws2.Range("G83").Value = 251674

$excel.Calculate()

# --- Update the selection / scroll position for each sheet view ---
$ws1.Activate()
$ws1.Range("C84").Select()
$excel.ActiveWindow.ScrollRow = 70

$ws2.Activate()
$ws2.Range("C84").Select()
$excel.ActiveWindow.ScrollRow = 68

$ws1.Activate()

$wb.Save()
